$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 01:27"

# Update country rows: label (where rank order changed) and statistics columns B:H

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 3354214
$ws.Range("C4").Value = 60287
$ws.Range("D4").Value = 1487935
$ws.Range("E4").Value = 1728893
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 715
$ws.Range("H4").Value = 137386

# Row 22: Colombia
$ws.Range("A22").Value = "Colombia"
$ws.Range("B22").Value = 145362
$ws.Range("C22").Value = 4586
$ws.Range("D22").Value = 61186
$ws.Range("E22").Value = 79057
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 194
$ws.Range("H22").Value = 5119

# Row 25: Argentina
$ws.Range("A25").Value = "Argentina"
$ws.Range("B25").Value = 97509
$ws.Range("C25").Value = 3449
$ws.Range("D25").Value = 41408
$ws.Range("E25").Value = 54291
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 1810

# Row 44: Panama
$ws.Range("A44").Value = "Panama"
$ws.Range("B44").Value = 44332
$ws.Range("C44").Value = 1075
$ws.Range("D44").Value = 22170
$ws.Range("E44").Value = 21269
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 30
$ws.Range("H44").Value = 893

# Row 52: Nigeria
$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 31987
$ws.Range("C52").Value = 664
$ws.Range("D52").Value = 13103
$ws.Range("E52").Value = 18160
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 15
$ws.Range("H52").Value = 724

# Row 53: Armenia
$ws.Range("A53").Value = "Armenia"
$ws.Range("B53").Value = 31392
$ws.Range("C53").Value = 489
$ws.Range("D53").Value = 19419
$ws.Range("E53").Value = 11414
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 13
$ws.Range("H53").Value = 559

# Row 54: Guatemala
$ws.Range("A54").Value = "Guatemala"
$ws.Range("B54").Value = 28598
$ws.Range("C54").Value = 979
$ws.Range("D54").Value = 4073
$ws.Range("E54").Value = 23353
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 33
$ws.Range("H54").Value = 1172

# Row 57: Ghana
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 24248
$ws.Range("C57").Value = 414
$ws.Range("D57").Value = 19831
$ws.Range("E57").Value = 4282
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 135

# Row 59: Japon
$ws.Range("A59").Value = "Japon"
$ws.Range("B59").Value = 21129
$ws.Range("C59").Value = 410
$ws.Range("D59").Value = 17849
$ws.Range("E59").Value = 2298
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 982

# Row 66: Camerun
$ws.Range("A66").Value = "Camerun"
$ws.Range("B66").Value = 15173
$ws.Range("C66").Value = 257
$ws.Range("D66").Value = 11928
$ws.Range("E66").Value = 2886
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 359

# Row 68: Chequia
$ws.Range("A68").Value = "Chequia"
$ws.Range("B68").Value = 13115
$ws.Range("C68").Value = 114
$ws.Range("D68").Value = 8227
$ws.Range("E68").Value = 4536
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 352

# Row 77: Noruega
$ws.Range("A77").Value = "Noruega"
$ws.Range("B77").Value = 8977
$ws.Range("C77").Value = 3
$ws.Range("D77").Value = 8138
$ws.Range("E77").Value = 587
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 252

# Row 82: Consejo Danes para los Refugiados
$ws.Range("A82").Value = "Consejo Danes para los Refugiados"
$ws.Range("B82").Value = 7971
$ws.Range("C82").Value = 66
$ws.Range("D82").Value = 3615
$ws.Range("E82").Value = 4167
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 189

# Row 90: Guinea
$ws.Range("A90").Value = "Guinea"
$ws.Range("B90").Value = 6044
$ws.Range("C90").Value = 75
$ws.Range("D90").Value = 4802
$ws.Range("E90").Value = 1205
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 37

# Row 107: Paraguay
$ws.Range("A107").Value = "Paraguay"
$ws.Range("B107").Value = 2820
$ws.Range("C107").Value = 84
$ws.Range("D107").Value = 1261
$ws.Range("E107").Value = 1538
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 21

# Row 130: Benin
$ws.Range("A130").Value = "Benin"
$ws.Range("B130").Value = 1378
$ws.Range("C130").Value = 93
$ws.Range("D130").Value = 557
$ws.Range("E130").Value = 795
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 3
$ws.Range("H130").Value = 26

# Row 131: Suazilandia
$ws.Range("A131").Value = "Suazilandia"
$ws.Range("B131").Value = 1311
$ws.Range("C131").Value = 54
$ws.Range("D131").Value = 656
$ws.Range("E131").Value = 637
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 18

# Row 132: Ruanda
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 1299
$ws.Range("C132").Value = 47
$ws.Range("D132").Value = 663
$ws.Range("E132").Value = 632
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 4

# Row 143: Uruguay
$ws.Range("A143").Value = "Uruguay"
$ws.Range("B143").Value = 986
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 896
$ws.Range("E143").Value = 60
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 30
